$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Plg"
$ws.Cells.Item(2,3).Value = "Plgrkt"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.465759
$ws.Cells.Item(2,8).Value = 1.397277
$ws.Cells.Item(2,9).Value = 0.9660838355812051
$ws.Cells.Item(2,10).Value = 0.9660838355812051
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 18.54216466666666
$ws.Cells.Item(2,14).Value = 55.62649399999999
$ws.Cells.Item(2,15).Value = 0.3703771243387364
$ws.Cells.Item(2,16).Value = 0.3703771243387364
$ws.Cells.Item(2,17).Value = 8.636180072982
$ws.Cells.Item(2,18).Value = 77.725620656838
$ws.Cells.Item(2,19).Value = 0.3578153528927034
$ws.Cells.Item(2,20).Value = 0.3578153528927034

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Plg"
$ws.Cells.Item(3,3).Value = "Plgrkt"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.465759
$ws.Cells.Item(3,8).Value = 1.397277
$ws.Cells.Item(3,9).Value = 0.9660838355812051
$ws.Cells.Item(3,10).Value = 0.9660838355812051
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 20.224881
$ws.Cells.Item(3,14).Value = 60.674643
$ws.Cells.Item(3,15).Value = 0.4039891457947978
$ws.Cells.Item(3,16).Value = 0.4039891457947978
$ws.Cells.Item(3,17).Value = 9.419920349679
$ws.Cells.Item(3,18).Value = 84.77928314711102
$ws.Cells.Item(3,19).Value = 0.390287383502613
$ws.Cells.Item(3,20).Value = 0.390287383502613

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Plg"
$ws.Cells.Item(4,3).Value = "Plgrkt"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.465759
$ws.Cells.Item(4,8).Value = 1.397277
$ws.Cells.Item(4,9).Value = 0.9660838355812051
$ws.Cells.Item(4,10).Value = 0.9660838355812051
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 11.295886
$ws.Cells.Item(4,14).Value = 33.887658
$ws.Cells.Item(4,15).Value = 0.2256337298664658
$ws.Cells.Item(4,16).Value = 0.2256337298664657
$ws.Cells.Item(4,17).Value = 5.261160567474001
$ws.Cells.Item(4,18).Value = 47.350445107266
$ws.Cells.Item(4,19).Value = 0.2179810991858888
$ws.Cells.Item(4,20).Value = 0.2179810991858887

# Row 5
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Plg"
$ws.Cells.Item(5,3).Value = "Plgrkt"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.01635133333333333
$ws.Cells.Item(5,8).Value = 0.049054
$ws.Cells.Item(5,9).Value = 0.03391616441879487
$ws.Cells.Item(5,10).Value = 0.03391616441879487
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 18.54216466666666
$ws.Cells.Item(5,14).Value = 55.62649399999999
$ws.Cells.Item(5,15).Value = 0.3703771243387364
$ws.Cells.Item(5,16).Value = 0.3703771243387364
$ws.Cells.Item(5,17).Value = 0.3031891151862222
$ws.Cells.Item(5,18).Value = 2.728702036676
$ws.Cells.Item(5,19).Value = 0.01256177144603301
$ws.Cells.Item(5,20).Value = 0.01256177144603301

# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Plg"
$ws.Cells.Item(6,3).Value = "Plgrkt"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.01635133333333333
$ws.Cells.Item(6,8).Value = 0.049054
$ws.Cells.Item(6,9).Value = 0.03391616441879487
$ws.Cells.Item(6,10).Value = 0.03391616441879487
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 20.224881
$ws.Cells.Item(6,14).Value = 60.674643
$ws.Cells.Item(6,15).Value = 0.4039891457947978
$ws.Cells.Item(6,16).Value = 0.4039891457947978
$ws.Cells.Item(6,17).Value = 0.330703770858
$ws.Cells.Item(6,18).Value = 2.976333937722
$ws.Cells.Item(6,19).Value = 0.01370176229218485
$ws.Cells.Item(6,20).Value = 0.01370176229218485

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Plg"
$ws.Cells.Item(7,3).Value = "Plgrkt"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.01635133333333333
$ws.Cells.Item(7,8).Value = 0.049054
$ws.Cells.Item(7,9).Value = 0.03391616441879487
$ws.Cells.Item(7,10).Value = 0.03391616441879487
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 11.295886
$ws.Cells.Item(7,14).Value = 33.887658
$ws.Cells.Item(7,15).Value = 0.2256337298664658
$ws.Cells.Item(7,16).Value = 0.2256337298664657
$ws.Cells.Item(7,17).Value = 0.1847027972813333
$ws.Cells.Item(7,18).Value = 1.662325175532
$ws.Cells.Item(7,19).Value = 0.007652630680576999
$ws.Cells.Item(7,20).Value = 0.007652630680576997

